$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 42) below the existing data (which ends at row 41).
# Force column A to be treated as literal text (not auto-parsed into a date
# serial number) so it matches the existing "yyyy/mm/dd"-as-text convention
# used by the rest of the sheet, then clear the resulting format so the cell
# keeps the sheet's default (unstyled) look like the other data rows.
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "2025/10/01"
$ws.Range("A42").ClearFormats()

$ws.Range("B42").Value = "水"
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 3
